# Applies the changes described by the commit "Opt with week and scenario runable".
#
# Sheets in the workbook:
#   bilateral_contract_data
#   electricity_demand
#   pool_price_scenarios          (not touched by this edit)
#   non_anticipativity_matrix
#   NU_non_anticipativity_matrix
#   pv_production

$wb = $excel.ActiveWorkbook

$wsBilateral = $wb.Worksheets.Item("bilateral_contract_data")
$wsDemand    = $wb.Worksheets.Item("electricity_demand")
$wsNonAnt    = $wb.Worksheets.Item("non_anticipativity_matrix")
$wsNuNonAnt  = $wb.Worksheets.Item("NU_non_anticipativity_matrix")
$wsPv        = $wb.Worksheets.Item("pv_production")

# ---------------------------------------------------------------------------
# 1) bilateral_contract_data: updated contract prices (column B) and new
#    custom column widths for C and D.
# ---------------------------------------------------------------------------
$wsBilateral.Range("B4").Value = 48.5
$wsBilateral.Range("B5").Value = 50
$wsBilateral.Range("B6").Value = 49.5
$wsBilateral.Range("B7").Value = 51
$wsBilateral.Range("B8").Value = 49
$wsBilateral.Range("B9").Value = 50

$wsBilateral.Columns.Item(3).ColumnWidth = 12.5
$wsBilateral.Columns.Item(4).ColumnWidth = 14.6

# ---------------------------------------------------------------------------
# 2) electricity_demand: row 2 (demand values) now follows the repeating
#    200 / 250 / 225 / 275 pattern all the way from column F to column CG.
# ---------------------------------------------------------------------------
$pattern = @(200, 250, 225, 275)
for ($col = 6; $col -le 85; $col++) {
    $wsDemand.Cells.Item(2, $col).Value = $pattern[($col - 2) % 4]
}

# ---------------------------------------------------------------------------
# 3) NU_non_anticipativity_matrix: cell X9 flips from 0 to 1.
# ---------------------------------------------------------------------------
$wsNuNonAnt.Range("X9").Value = 1

# ---------------------------------------------------------------------------
# 4) Sheet view / selection updates (this also drives which sheet ends up
#    as the active/selected tab once the workbook is saved).
# ---------------------------------------------------------------------------
$wsDemand.Activate()
$wsDemand.Range("O5").Select()

$wsNonAnt.Activate()
$wsNonAnt.Range("G23").Select()

$wsNuNonAnt.Activate()
$wsNuNonAnt.Range("AC11").Select()

$wsPv.Activate()
$wsPv.Range("K22").Select()

$wsBilateral.Activate()
$wsBilateral.Range("D10").Select()
